$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Apply-Runs($range, $runs) {
    # $runs is an array of hashtables: @{ Text = "..."; Bold = $true/$false; Italic = $true/$false }
    $full = ""
    foreach ($r in $runs) { $full += $r.Text }
    $range.Value = $full
    $pos = 1
    foreach ($r in $runs) {
        $len = $r.Text.Length
        if ($len -gt 0) {
            $seg = $range.Characters($pos, $len)
            $seg.Font.Bold = [bool]$r.Bold
            $seg.Font.Italic = [bool]$r.Italic
        }
        $pos += $len
    }
}

# --- 1) Update the "Last updated" banner in B1 ---
$ws.Range("B1").Value = "Last updated: 02/15/2017 12:28 PM by Carlos"

# --- 2) Update the date inside the rich-text cell E6 (2/5 -> 2/15), keep the rest unchanged ---
$e6Runs = @(
    @{ Text = "Please complete the following by 2/15"; Bold = $true;  Italic = $true  }
    @{ Text = ":";                                      Bold = $true;  Italic = $false }
    @{ Text = " ";                                       Bold = $false; Italic = $false }
    @{ Text = "Carlos:";                                 Bold = $true;  Italic = $false }
    @{ Text = " Update Project Milestones; ";            Bold = $false; Italic = $false }
    @{ Text = "Courtnie: ";                               Bold = $true;  Italic = $false }
    @{ Text = "Update House of Quality; ";                Bold = $false; Italic = $false }
    @{ Text = "Lucas:";                                   Bold = $true;  Italic = $false }
    @{ Text = " goals and objectives;  ";                 Bold = $false; Italic = $false }
    @{ Text = "Patrick:";                                 Bold = $true;  Italic = $false }
    @{ Text = " Project Block Diagram; ";                 Bold = $false; Italic = $false }
    @{ Text = "Complete by 2/22";                         Bold = $true;  Italic = $true  }
    @{ Text = ": ";                                       Bold = $false; Italic = $false }
    @{ Text = "Carlos";                                   Bold = $true;  Italic = $false }
    @{ Text = ": Smoke Sensors Research & Docs; ";        Bold = $false; Italic = $false }
    @{ Text = "Courtnie";                                 Bold = $true;  Italic = $false }
    @{ Text = ": Battery Research & Docs; ";              Bold = $false; Italic = $false }
    @{ Text = "Lucas:";                                   Bold = $true;  Italic = $false }
    @{ Text = " Microcontroller Research & Docs; ";       Bold = $false; Italic = $false }
    @{ Text = "Patrick";                                  Bold = $true;  Italic = $false }
    @{ Text = ": Wireless Communication Research & Docs;"; Bold = $false; Italic = $false }
)
Apply-Runs $ws.Range("E6") $e6Runs

Write-Output "done part 1"

# --- 3) Add the new meeting row (row 7, 2/15/2017 meeting) ---

# Copy formatting from row 6 so the new row gets identical styles
# (date format + centered in col A, wrap text in D/E/F).
$ws.Range("A6:F6").Copy()
$ws.Range("A7:F7").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Rows(7).RowHeight = 90

$ws.Range("A7").Value = 42781
$ws.Range("B7").Value = "Discuss Updates for Initial Document/Proposal. Make sure everyone has Research Assigned, and Understands what is Required."
$ws.Range("C7").Value = "Carlos, Courtnie, Lucas, Patrick"
$ws.Range("D7").Value = 'Assigned Research for everyone. Will submit "Updates" documented by 2/16.'
$ws.Range("F7").Value = "Will look into sponsors later on in project. Add Software Flowchart to document whenever possible."

# E7 is rich text: bold+italic lead-in, then alternating bold names / plain descriptions.
$e7Runs = @(
    @{ Text = "Please complete Research on at least one of the following by 2/22"; Bold = $true;  Italic = $true  }
    @{ Text = ": Carlos";                                                          Bold = $true;  Italic = $false }
    @{ Text = ": Smoke Sensors, Direction/Location Algorithm Research & Docs; ";   Bold = $false; Italic = $false }
    @{ Text = "Courtnie";                                                          Bold = $true;  Italic = $false }
    @{ Text = ": Battery/Power Monitoring, Alarm System Components Research & Docs; "; Bold = $false; Italic = $false }
    @{ Text = "Lucas:";                                                            Bold = $true;  Italic = $false }
    @{ Text = " Microcontroller and Processor Research & Docs; ";                  Bold = $false; Italic = $false }
    @{ Text = "Patrick";                                                           Bold = $true;  Italic = $false }
    @{ Text = ": Wireless Communication Research and Processor & Docs;";           Bold = $false; Italic = $false }
)
Apply-Runs $ws.Range("E7") $e7Runs

Write-Output "done part 2"

